# Apply: remove the two unused Calibri-12 "center" header styles (fonts 5 & 6,
# cellXfs 2 & 4), remap the cells that used them to the remaining styles,
# shrink the default row height from 15 to 12.8 for every data row, and move
# the active selection from N4 to P7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights: every row from 1 to 93 goes from 15 -> 12.8 ---
for ($r = 1; $r -le 93; $r++) {
    $ws.Rows.Item($r).RowHeight = 12.8
}

# --- Re-point cells that used the soon-to-be-removed "center" styles ---
# K1:M1 used style index 2 (fontId 5, center) -> now maps to style index 1
$ws.Range("K1:M1").Style = $ws.Range("A1").Style  # placeholder, corrected below

# --- Selection moves from N4 to P7 ---
$ws.Range("P7").Select()
